$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename topic prefixes from "DSMR-API/" to "TOP TOPIC/" in column A (rows 2-27)
for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Text
    if ($val -ne $null -and $val.ToString().StartsWith("DSMR-API/")) {
        $cell.Value = $val.ToString().Replace("DSMR-API/", "TOP TOPIC/")
    }
}

# Widen column B (target stored width 65.33203125; ColumnWidth input is offset by ~5/6 by the host)
$ws.Columns.Item(2).ColumnWidth = 64.49869791666667

# Change selection to A3
$ws.Range("A3").Select()
